# Generate Report for Archive
# - Update status from "Ready for handoff" to "In Translation" everywhere it
#   appears (Overview!E2:E3/F2:F3 and the Status column (C2:C3) on the
#   per-locale sheets), then resize the now-narrower Status columns to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: "zh-cn" (col E) and "de-de" (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Narrow the two status columns now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: "Status" column (col C) ---
$localeSheets = @("zh-cn", "de-de")
foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
